$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): extend sequence with P1=14, Q1=15, matching O1's style ---
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# --- Data rows 2-25: add new columns P and Q (value 2), and swap I<->K and M<->O ---
for ($r = 2; $r -le 25; $r++) {
    $iVal = $ws.Cells.Item($r, 9).Value()   # column I
    $kVal = $ws.Cells.Item($r, 11).Value()  # column K
    $mVal = $ws.Cells.Item($r, 13).Value()  # column M
    $oVal = $ws.Cells.Item($r, 15).Value()  # column O

    $ws.Cells.Item($r, 9).Value = $kVal   # I <- old K
    $ws.Cells.Item($r, 11).Value = $iVal  # K <- old I
    $ws.Cells.Item($r, 13).Value = $oVal  # M <- old O
    $ws.Cells.Item($r, 15).Value = $mVal  # O <- old M

    $ws.Cells.Item($r, 16).Value = 2  # P
    $ws.Cells.Item($r, 17).Value = 2  # Q
}
